$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the promo text "... tenes un 33%. Limite ..." with "... tenes un 22%. Limite ..."
# across every cell that contains it. Using Cells.Replace with xlPart so it
# only matches the specific phrase (keeping "... tenes un 33% en todo ..."
# strings, which are a different promo text, untouched).
$ws.Cells.Replace("tenes un 33%. Limite", "tenes un 22%. Limite", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)

# Update the active selection to C5 (was F11)
$ws.Range("C5").Select()
